{"js": "// fix(docx): fix OOXMLValidator error on *Tok output\n//\n// The <w:rPr> of several \"Tok\" character styles (used for syntax-\n// highlighted source code) serialized <w:color/> before <w:b/>/<w:i/>,\n// which violates the CT_RPr sequence in wml.xsd (toggle properties\n// such as b/i must precede color). Re-assigning each affected style's\n// Font.bold / Font.italic property (to its own current value) forces\n// the host to re-serialize <w:rPr>'s children in schema order without\n// changing any formatting.\n\n// styleId -> which toggle propertie(s) to \"touch\" (re-set to their own\n// current value) so the rPr gets rewritten in schema order. Only the\n// properties already present on the style are touched, so no new\n// <w:b w:val=\"0\"/> / <w:i w:val=\"0\"/> elements get introduced.\nconst styleProps = {\n  KeywordTok: [\"bold\"],\n  ImportTok: [\"bold\"],\n  CommentTok: [\"italic\"],\n  DocumentationTok: [\"italic\"],\n  AnnotationTok: [\"bold\", \"italic\"],\n  CommentVarTok: [\"bold\", \"italic\"],\n  ControlFlowTok: [\"bold\"],\n  InformationTok: [\"bold\", \"italic\"],\n  WarningTok: [\"bold\", \"italic\"],\n  AlertTok: [\"bold\"],\n  ErrorTok: [\"bold\"],\n};\n\nconst styleNames = Object.keys(styleProps);\nconst fonts = {};\nfor (const name of styleNames) {\n  const style = context.document.getStyles().getByName(name);\n  const font = style.font;\n  font.load([\"bold\", \"italic\"]);\n  fonts[name] = font;\n}\n\nawait context.sync();\n\nfor (const name of styleNames) {\n  const font = fonts[name];\n  for (const prop of styleProps[name]) {\n    if (prop === \"bold\") {\n      font.bold = font.bold;\n    } else if (prop === \"italic\") {\n      font.italic = font.italic;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# fix(docx): fix OOXMLValidator error on *Tok character styles\n#\n# The <w:rPr> of several \"Tok\" character styles (used for syntax-\n# highlighted source code) had <w:color/> written before <w:b/>/<w:i/>,\n# which violates the CT_RPr sequence in wml.xsd (toggle properties like\n# b/i must precede color). Re-assigning the Font.Bold / Font.Italic\n# properties on each affected style forces Word to re-serialize the\n# <w:rPr> children in schema order, without changing any values.\n\n$d = $word.ActiveDocument\n\n# styleId -> which toggle propertie(s) to \"touch\" (re-set to their own\n# current value) so the rPr gets rewritten in schema order. Only the\n# properties already present on the style are touched, so no new\n# <w:b w:val=\"0\"/> / <w:i w:val=\"0\"/> elements get introduced.\n$styleProps = @{\n    \"KeywordTok\"       = @(\"Bold\")\n    \"ImportTok\"        = @(\"Bold\")\n    \"CommentTok\"       = @(\"Italic\")\n    \"DocumentationTok\" = @(\"Italic\")\n    \"AnnotationTok\"    = @(\"Bold\", \"Italic\")\n    \"CommentVarTok\"    = @(\"Bold\", \"Italic\")\n    \"ControlFlowTok\"   = @(\"Bold\")\n    \"InformationTok\"   = @(\"Bold\", \"Italic\")\n    \"WarningTok\"       = @(\"Bold\", \"Italic\")\n    \"AlertTok\"         = @(\"Bold\")\n    \"ErrorTok\"         = @(\"Bold\")\n}\n\nforeach ($styleName in $styleProps.Keys) {\n    $style = $d.Styles($styleName)\n    $font = $style.Font\n    foreach ($prop in $styleProps[$styleName]) {\n        if ($prop -eq \"Bold\") {\n            $font.Bold = $font.Bold\n        } elseif ($prop -eq \"Italic\") {\n            $font.Italic = $font.Italic\n        }\n    }\n}\n"}
